$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = "93_referral_statement"
$ws.Range("F43").Value = "ppe || 18_hazards_to_humans_and_domestic_animals"
$ws.Range("F45").Value = "ppe"
$ws.Range("F46").Value = "ppe"
$ws.Range("F49").Value = "application instructions || env warning - species"
$ws.Range("F50").Value = "application instructions || env warning - species"
$ws.Range("F51").Value = "application instructions || env warning - species"
$ws.Range("F52").Value = "env warning - water"
$ws.Range("F53").Value = "env warning - water || off target movement"
$ws.Range("F56").Value = "application instructions"
$ws.Range("F60").Value = "application instructions"
$ws.Range("F61").Value = "application instructions"
$ws.Range("F62").Value = "pollinator"
$ws.Range("F63").Value = "135_product_information"
$ws.Range("F67").Value = "use restrictions"
$ws.Range("F70").Value = "mixing"
$ws.Range("F71").Value = "mixing"
$ws.Range("F73").Value = "mixing"
$ws.Range("F74").Value = "use restrictions"
$ws.Range("F75").Value = "use restrictions"
$ws.Range("F88").Value = "application instructions"
$ws.Range("F100").Value = "off target movement"
$ws.Range("F101").Value = "off target movement"
$ws.Range("F102").Value = "off target movement"
$ws.Range("F104").Value = "off target movement"
$ws.Range("F105").Value = "off target movement"
$ws.Range("F108").Value = "off target movement"
$ws.Range("F109").Value = "off target movement"
$ws.Range("F112").Value = "use restrictions || application instructions"
$ws.Range("F113").Value = "use restrictions || application instructions"
$ws.Range("F115").Value = "application instructions"
$ws.Range("F116").Value = "safety procedures || application instructions"
$ws.Range("F117").Value = "mixing"
$ws.Range("F118").Value = "mixing"
$ws.Range("F119").Value = "use restrictions"
$ws.Range("F120").Value = "use restrictions || mixing"
$ws.Range("F121").Value = "application instructions"
$ws.Range("F123").Value = "application instructions"
$ws.Range("F124").Value = "use restrictions || application instructions"
$ws.Range("F125").Value = "use restrictions || application instructions"
$ws.Range("F126").Value = "use restrictions"
$ws.Range("F127").Value = "irrigation || application instructions || chemigation"
$ws.Range("F128").Value = "safety procedures"
$ws.Range("F130").Value = "safety procedures"
$ws.Range("F131").Value = "irrigation"
$ws.Range("F132").Value = "irrigation"
$ws.Range("F133").Value = "irrigation"
$ws.Range("F135").Value = "irrigation"
$ws.Range("F582").Value = "mixing"
$ws.Range("F584").Value = "use restrictions"
